# Updates the latest crypto price-snapshot columns (Price / Volume(1h))
# on the active worksheet to match the newly scraped symbol list.
# Values are written as plain text (matching the source data, which
# stores "304.89" / "-0.79%" etc. as literal strings, not numbers),
# so NumberFormat is forced to Text ("@") for the write and the
# cell's original Style is restored immediately afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value
$updates = @{
    "D2" = "304.89"
    "E2" = "-0.79%"
    "D3" = "35.90"
    "E3" = "-1.61%"
    "D4" = "5.008"
    "E4" = "-1.67%"
    "D5" = "0.08059"
    "E5" = "-0.62%"
    "D6" = "1.886"
    "E6" = "-4.51%"
    "E7" = "0.69%"
    "D8" = "7.849"
    "E8" = "1.11%"
    "D9" = "0.9325"
    "E9" = "-0.32%"
    "D10" = "0.1312"
    "E10" = "-9.62%"
    "D11" = "0.1908"
    "E11" = "-1.21%"
    "D12" = "0.09228"
    "E12" = "0.92%"
    "D13" = "0.03516"
    "E13" = "-0.73%"
    "D14" = "0.09890"
    "E14" = "1.01%"
    "D15" = "0.001417"
    "E15" = "-1.02%"
    "D16" = "0.006713"
    "E16" = "15.56%"
    "D17" = "3.609"
    "E17" = "2.43%"
    "D18" = "3.139"
    "E18" = "5.27%"
    "E19" = "0.55%"
    "D20" = "0.1345"
    "E20" = "3.29%"
    "D21" = "5.248"
    "E21" = "5.90%"
    "D22" = "0.2534"
    "E22" = "5.50%"
    "D23" = "0.04429"
    "E23" = "-1.68%"
    "E24" = "1.93%"
    "D25" = "0.004714"
    "E25" = "-2.98%"
    "D26" = "0.0001302"
    "E26" = "4.99%"
    "E27" = "-29.58%"
    "D39" = "0.01941"
    "E39" = "-1.75%"
    "D40" = "0.05174"
    "E40" = "6.19%"
    "D41" = "0.007548"
    "E41" = "-0.10%"
    "D42" = "0.01019"
    "E42" = "-7.77%"
    "D43" = "0.1373"
    "E43" = "0.47%"
    "D44" = "0.002164"
    "E44" = "2.51%"
    "D45" = "0.01080"
    "E45" = "10.60%"
    "D46" = "0.00006366"
    "E46" = "-0.14%"
    "E47" = "0.14%"
    "D48" = "65.22"
    "E48" = "0.85%"
    "D49" = "0.001661"
    "E49" = "39.44%"
    "E50" = "0.14%"
    "E51" = "0.14%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = $origStyle
}
